$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '60.259.19'
Set-TextCell 2 5 '  -2.82%  '

Set-TextCell 3 4 '3.304.12'
Set-TextCell 3 5 '  -3.27%  '

Set-TextCell 4 4 '1.00'
Set-TextCell 4 5 '  +0.10%  '

Set-TextCell 5 4 '558.04'
Set-TextCell 5 5 '  -3.15%  '

Set-TextCell 6 4 '142.11'
Set-TextCell 6 5 '  -4.12%  '

Set-TextCell 7 5 '  +0.06%  '

Set-TextCell 8 4 '3.305.83'
Set-TextCell 8 5 '  -3.19%  '

Set-TextCell 9 5 '  -2.48%  '

Set-TextCell 10 4 '7.85'
Set-TextCell 10 5 '  -1.46%  '

Set-TextCell 11 5 '  -3.81%  '

Set-TextCell 12 5 '  -2.08%  '

Set-TextCell 13 4 '3.875.16'
Set-TextCell 13 5 '  -2.96%  '

Set-TextCell 14 5 '  +0.44%  '

Set-TextCell 15 4 '26.83'
Set-TextCell 15 5 '  -5.53%  '

Set-TextCell 16 4 '3.298.12'
Set-TextCell 16 5 '  -3.15%  '

Set-TextCell 17 5 '  -3.32%  '

Set-TextCell 18 4 '60.263.55'
Set-TextCell 18 5 '  -2.77%  '

Set-TextCell 19 4 '6.19'
Set-TextCell 19 5 '  -3.49%  '

Set-TextCell 20 4 '14.43'
Set-TextCell 20 5 '  -0.89%  '

Set-TextCell 21 4 '8.63'
Set-TextCell 21 5 '  -3.75%  '

Set-TextCell 22 4 '374.42'
Set-TextCell 22 5 '  -1.71%  '

Set-TextCell 23 4 '74.24'
Set-TextCell 23 5 '  -0.81%  '

Set-TextCell 24 5 '  -4.34%  '

Set-TextCell 25 5 '  +0.02%  '

Set-TextCell 26 4 '3.446.74'
Set-TextCell 26 5 '  -3.41%  '

Set-TextCell 27 5 '  -7.91%  '

Set-TextCell 28 4 '0.171'
Set-TextCell 28 5 '  -4.55%  '

Set-TextCell 29 5 '  +0.04%  '

Set-TextCell 30 4 '7.24'
Set-TextCell 30 5 '  -4.80%  '

Set-TextCell 31 5 '  -0.09%  '

Set-TextCell 32 4 '7.62'
Set-TextCell 32 5 '  -3.78%  '

Set-TextCell 33 4 '2.04'
Set-TextCell 33 5 '  -3.69%  '

Set-TextCell 34 4 '22.55'
Set-TextCell 34 5 '  -2.25%  '

Set-TextCell 35 5 '  -5.39%  '

Set-TextCell 36 5 '  -6.24%  '

Set-TextCell 37 4 '165.88'
Set-TextCell 37 5 '  -2.27%  '

Set-TextCell 38 5 '  -5.14%  '

Set-TextCell 39 5 '  -2.83%  '

Set-TextCell 40 2 'RenzoRestakedETH'
Set-TextCell 40 3 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextCell 40 4 '3.336.32'
Set-TextCell 40 5 '  -3.12%  '

Set-TextCell 41 2 'EnergySwap'
Set-TextCell 41 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 41 4 '26.76'
Set-TextCell 41 5 '  -11.40%  '

Set-TextCell 42 5 '  -5.79%  '

Set-TextCell 43 4 '41.98'
Set-TextCell 43 5 '  -1.05%  '

Set-TextCell 44 4 '0.752'
Set-TextCell 44 5 '  -3.27%  '

Set-TextCell 45 5 '  -3.96%  '

Set-TextCell 46 5 '  -5.16%  '

Set-TextCell 47 4 '1.11'
Set-TextCell 47 5 '  -4.61%  '

Set-TextCell 48 4 '2.373.76'
Set-TextCell 48 5 '  -6.61%  '

Set-TextCell 49 5 '  +0.17%  '

Set-TextCell 50 5 '  -5.86%  '

Set-TextCell 51 5 '  -6.48%  '
